$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B9: status changes from "offen" to "in Arbeit" (Mails will be sent now after order)
$ws.Range("B9").Value = "in Arbeit"
$ws.Range("B9").Style = "Neutral"

# Update the active selection to D9
$ws.Range("D9").Select()
